$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13 (hunk 0)
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""

# Row 31 (hunk 1)
$ws.Range("H31").Value = 3023.7144
$ws.Range("I31").Value = 833.2
$ws.Range("J31").Value = 8500
$ws.Range("K31").Value = 2499.6
$ws.Range("L31").Value = 25500
$ws.Range("M31").Value = -2269.6
$ws.Range("N31").Value = -25960

# Row 41 (hunk 2)
$ws.Range("H41").Value = 2431
$ws.Range("I41").Value = 3639
$ws.Range("J41").Value = 981.4
$ws.Range("K41").Value = 3639
$ws.Range("L41").Value = 981.4
$ws.Range("M41").Value = -3199
$ws.Range("N41").Value = -1861.4

# Row 98 (hunk 3)
$ws.Range("H98").Value = 1220
$ws.Range("I98").Value = 1220
$ws.Range("K98").Value = 1220
$ws.Range("M98").Value = 278

# Row 122 (hunk 4)
$ws.Range("H122").Value = 1220
$ws.Range("I122").Value = 1220
$ws.Range("K122").Value = 3660
$ws.Range("M122").Value = -1210

$ws = $wb.Worksheets.Item("ARM")
# Row 13 (hunk 5)
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""

# Row 32 (hunk 6)
$ws.Range("H32").Value = 4709
$ws.Range("I32").Value = 5100
$ws.Range("K32").Value = 5100
$ws.Range("M32").Value = -4813

# Row 88 (hunk 7)
$ws.Range("H88").Value = 737.5
$ws.Range("I88").Value = 450
$ws.Range("K88").Value = 450
$ws.Range("M88").Value = -44

# Row 91 (hunk 8)
$ws.Range("H91").Value = 737.5
$ws.Range("I91").Value = 450
$ws.Range("K91").Value = 450
$ws.Range("M91").Value = 954

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (hunk 9)
$ws.Range("H86").Value = 1803.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1803.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1803.5
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -4049.5

# Row 89 (hunk 10)
$ws.Range("H89").Value = 1803.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1803.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 9017.5
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -20249.5

# Row 94 (hunk 11)
$ws.Range("H94").Value = 1999
$ws.Range("I94").Value = 1999
$ws.Range("K94").Value = 1999
$ws.Range("M94").Value = -1548

# Row 124 (hunk 12)
$ws.Range("H124").Value = 47865.8
$ws.Range("J124").Value = 47865.8
$ws.Range("L124").Value = 47865.8
$ws.Range("N124").Value = -57685.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 13)
$ws.Range("H16").Value = 1490
$ws.Range("I16").Value = 1470
$ws.Range("K16").Value = 1470
$ws.Range("M16").Value = -1183

# Row 20 (hunk 14)
$ws.Range("H20").Value = 147000
$ws.Range("J20").Value = 147000
$ws.Range("L20").Value = 147000
$ws.Range("N20").Value = -147472

# Row 30 (hunk 15)
$ws.Range("H30").Value = 147000
$ws.Range("J30").Value = 147000
$ws.Range("L30").Value = 147000
$ws.Range("N30").Value = -147182

# Row 31 (hunk 16)
$ws.Range("H31").Value = 2679
$ws.Range("I31").Value = 2723.75
$ws.Range("K31").Value = 2723.75
$ws.Range("M31").Value = -2428.75

# Row 34 (hunk 17)
$ws.Range("H34").Value = 2679
$ws.Range("I34").Value = 2723.75
$ws.Range("K34").Value = 2723.75
$ws.Range("M34").Value = -2521.75

# Row 113 (hunk 18)
$ws.Range("H113").Value = 1490
$ws.Range("I113").Value = 1470
$ws.Range("K113").Value = 1470
$ws.Range("M113").Value = 700

# Row 122 (hunk 19)
$ws.Range("H122").Value = 21598.4
$ws.Range("I122").Value = 14498
$ws.Range("J122").Value = 50000
$ws.Range("K122").Value = 43494
$ws.Range("L122").Value = 150000
$ws.Range("M122").Value = -41044
$ws.Range("N122").Value = -154900

# Row 128 (hunk 20)
$ws.Range("H128").Value = 147000
$ws.Range("J128").Value = 147000
$ws.Range("L128").Value = 147000
$ws.Range("N128").Value = -156960

# Row 134 (hunk 21)
$ws.Range("H134").Value = 2931.6875
$ws.Range("I134").Value = 1613.5
$ws.Range("K134").Value = 4840.5
$ws.Range("M134").Value = -2305.5

$ws = $wb.Worksheets.Item("CUL")
# Row 11 (hunk 22)
$ws.Range("H11").Value = 418.57144
$ws.Range("I11").Value = 438.33334
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 1315.00002
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -1175.00002
$ws.Range("N11").Value = -1180

# Row 34 (hunk 23)
$ws.Range("H34").Value = 4244.4443
$ws.Range("I34").Value = 500
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = -1416

# Row 133 (hunk 24)
$ws.Range("H133").Value = 15000
$ws.Range("J133").Value = 15000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120

$ws = $wb.Worksheets.Item("GSM")
# Row 98 (hunk 25)
$ws.Range("H98").Value = 20778
$ws.Range("J98").Value = 20778
$ws.Range("L98").Value = 20778
$ws.Range("N98").Value = -26768

# Row 122 (hunk 26)
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

# Row 126 (hunk 27)
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 28)
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1205

# Row 27 (hunk 29)
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 1500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1393

# Row 40 (hunk 30)
$ws.Range("H40").Value = 6618
$ws.Range("I40").Value = 6618
$ws.Range("K40").Value = 6618
$ws.Range("M40").Value = -6482

# Row 93 (hunk 31)
$ws.Range("H93").Value = 4998.25
$ws.Range("I93").Value = 4998.25
$ws.Range("K93").Value = 4998.25
$ws.Range("M93").Value = -3750.25

# Row 136 (hunk 32)
$ws.Range("H136").Value = 14165.667
$ws.Range("I136").Value = 6998.8
$ws.Range("K136").Value = 20996.4
$ws.Range("M136").Value = -18446.4

$ws = $wb.Worksheets.Item("WVR")
# Row 5 (hunk 33)
$ws.Range("H5").Value = 12502525
$ws.Range("I5").Value = 3367
$ws.Range("J5").Value = 50000000
$ws.Range("K5").Value = 3367
$ws.Range("L5").Value = 50000000
$ws.Range("M5").Value = -3255
$ws.Range("N5").Value = -50000224

# Row 13 (hunk 34)
$ws.Range("H13").Value = 225
$ws.Range("I13").Value = 225
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 225
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -85
$ws.Range("N13").Value = ""

# Row 96 (hunk 35)
$ws.Range("H96").Value = 1225
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1300
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1300
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4046

# Row 132 (hunk 36)
$ws.Range("H132").Value = 4141.615
$ws.Range("I132").Value = 4344.6
$ws.Range("J132").Value = 3465
$ws.Range("K132").Value = 13033.8
$ws.Range("L132").Value = 10395
$ws.Range("M132").Value = -10503.8
$ws.Range("N132").Value = -15455

# Row 137 (hunk 37)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""
